$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the "to Mars ... Europa" paragraph: the wormhole destination and
#    meteoroid/ship-repair / mission-ending sentences were reworded.
# ---------------------------------------------------------------------------
$oldTail = " to Mars and from there he has to travel across 100 asteroids, meanwhile he has several random chances of being teleported by wormholes to different places in our solar system or even in an unknown area of the universe. He also has a random chance of getting hit by meteoroids on his way to Europa."
$newTail = " to Mars and from there he has to travel across 100 asteroids, meanwhile he has several random chances of being teleported by wormholes to the Sun, an unknown area of the universe, or get hit by a meteoroid. If so, he’ll have to get back to Earth to repair his ship. At the end of the day his mission is to make it to Europa in one piece!"

$rng = $d.Content
$rng.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

Write-Output "done"
